# Updates the cryptocurrency price/volume table with freshly scraped values.
# Generated from the authoritative row-by-row change list below; each entry is
# applied as $ws.Range(<cell>).Value = <new text>. Column D ("Price") is stored
# as text in this sheet (e.g. "67.829.49", "1.00", "0.550"), so those cells are
# first marked with a Text number format -- otherwise Excel would silently
# reinterpret a numeric-looking string as a Number and drop formatting such as
# trailing zeros or the dotted-thousands style used for large prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '67.829.49'; ForceText = $true },
    @{ Cell = 'E2'; Value = '  +1.55%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '2.617.06'; ForceText = $true },
    @{ Cell = 'E3'; Value = '  +1.30%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '601.32'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  +1.39%  '; ForceText = $false },
    @{ Cell = 'D6'; Value = '154.38'; ForceText = $true },
    @{ Cell = 'E7'; Value = '  -0.03%  '; ForceText = $false },
    @{ Cell = 'D8'; Value = '0.550'; ForceText = $true },
    @{ Cell = 'E8'; Value = '  +1.47%  '; ForceText = $false },
    @{ Cell = 'D9'; Value = '2.616.13'; ForceText = $true },
    @{ Cell = 'E9'; Value = '  +1.30%  '; ForceText = $false },
    @{ Cell = 'E10'; Value = '  +10.55%  '; ForceText = $false },
    @{ Cell = 'E11'; Value = '  +0.90%  '; ForceText = $false },
    @{ Cell = 'E12'; Value = '  +0.93%  '; ForceText = $false },
    @{ Cell = 'E13'; Value = '  -0.85%  '; ForceText = $false },
    @{ Cell = 'D14'; Value = '27.61'; ForceText = $true },
    @{ Cell = 'E14'; Value = '  -2.06%  '; ForceText = $false },
    @{ Cell = 'E15'; Value = '  +3.50%  '; ForceText = $false },
    @{ Cell = 'D16'; Value = '3.095.17'; ForceText = $true },
    @{ Cell = 'E16'; Value = '  +1.36%  '; ForceText = $false },
    @{ Cell = 'D17'; Value = '67.781.07'; ForceText = $true },
    @{ Cell = 'E17'; Value = '  +1.46%  '; ForceText = $false },
    @{ Cell = 'D18'; Value = '2.617.78'; ForceText = $true },
    @{ Cell = 'E18'; Value = '  +1.41%  '; ForceText = $false },
    @{ Cell = 'D19'; Value = '11.19'; ForceText = $true },
    @{ Cell = 'E19'; Value = '  -0.64%  '; ForceText = $false },
    @{ Cell = 'D20'; Value = '365.90'; ForceText = $true },
    @{ Cell = 'E20'; Value = '  +3.37%  '; ForceText = $false },
    @{ Cell = 'E21'; Value = '  -1.59%  '; ForceText = $false },
    @{ Cell = 'E22'; Value = '  -0.49%  '; ForceText = $false },
    @{ Cell = 'E23'; Value = '  -1.98%  '; ForceText = $false },
    @{ Cell = 'D24'; Value = '1.00'; ForceText = $true },
    @{ Cell = 'E24'; Value = '  -0.06%  '; ForceText = $false },
    @{ Cell = 'D25'; Value = '70.31'; ForceText = $true },
    @{ Cell = 'E25'; Value = '  +4.60%  '; ForceText = $false },
    @{ Cell = 'D26'; Value = '9.85'; ForceText = $true },
    @{ Cell = 'E26'; Value = '  -6.34%  '; ForceText = $false },
    @{ Cell = 'E27'; Value = '  +0.85%  '; ForceText = $false },
    @{ Cell = 'E28'; Value = '  +1.18%  '; ForceText = $false },
    @{ Cell = 'D29'; Value = '576.37'; ForceText = $true },
    @{ Cell = 'E29'; Value = '  -3.96%  '; ForceText = $false },
    @{ Cell = 'E30'; Value = '  +0.08%  '; ForceText = $false },
    @{ Cell = 'E31'; Value = '  -2.19%  '; ForceText = $false },
    @{ Cell = 'D32'; Value = '7.90'; ForceText = $true },
    @{ Cell = 'E32'; Value = '  -1.99%  '; ForceText = $false },
    @{ Cell = 'E33'; Value = '  +0.61%  '; ForceText = $false },
    @{ Cell = 'E34'; Value = '  -1.49%  '; ForceText = $false },
    @{ Cell = 'E35'; Value = '  +0.08%  '; ForceText = $false },
    @{ Cell = 'E36'; Value = '  -3.08%  '; ForceText = $false },
    @{ Cell = 'D37'; Value = '4.92'; ForceText = $true },
    @{ Cell = 'E37'; Value = '  -1.74%  '; ForceText = $false },
    @{ Cell = 'D38'; Value = '158.62'; ForceText = $true },
    @{ Cell = 'E38'; Value = '  +3.12%  '; ForceText = $false },
    @{ Cell = 'D39'; Value = '19.34'; ForceText = $true },
    @{ Cell = 'E39'; Value = '  +0.78%  '; ForceText = $false },
    @{ Cell = 'E40'; Value = '  +0.33%  '; ForceText = $false },
    @{ Cell = 'E41'; Value = '  +3.44%  '; ForceText = $false },
    @{ Cell = 'D42'; Value = '5.35'; ForceText = $true },
    @{ Cell = 'E42'; Value = '  -1.78%  '; ForceText = $false },
    @{ Cell = 'E43'; Value = '  -1.70%  '; ForceText = $false },
    @{ Cell = 'E44'; Value = '  -0.83%  '; ForceText = $false },
    @{ Cell = 'E45'; Value = '  +0.11%  '; ForceText = $false },
    @{ Cell = 'D46'; Value = '16.42'; ForceText = $true },
    @{ Cell = 'E46'; Value = '  -0.13%  '; ForceText = $false },
    @{ Cell = 'D47'; Value = '156.93'; ForceText = $true },
    @{ Cell = 'E47'; Value = '  +0.91%  '; ForceText = $false },
    @{ Cell = 'E48'; Value = '  -7.17%  '; ForceText = $false },
    @{ Cell = 'E49'; Value = '  +0.09%  '; ForceText = $false },
    @{ Cell = 'D50'; Value = '20.92'; ForceText = $true },
    @{ Cell = 'E50'; Value = '  -1.85%  '; ForceText = $false },
    @{ Cell = 'B51'; Value = 'Mantle'; ForceText = $false },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; ForceText = $false },
    @{ Cell = 'D51'; Value = '0.623'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  +1.64%  '; ForceText = $false }
)

foreach ($u in $updates) {
    if ($u.ForceText) {
        # Pre-format as Text so a numeric-looking string (e.g. "1.00", "0.550")
        # is stored verbatim instead of being normalized into a Number.
        $ws.Range($u.Cell).NumberFormat = "@"
    }
    $ws.Range($u.Cell).Value = $u.Value
}
